# Updated Master data as per 16th May Refresh
# Append 3 new rows (34-36) to the reg_center_user_machine_h master data sheet,
# following the same pattern as the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows: regcntr_id, usr_id, machine_id, lang_code, is_active, cr_by, cr_dtimes, eff_dtimes
$newRows = @(
    @(10005, 110033, 10005, "eng", $true, "superadmin", "now()", "now()"),
    @(10005, 110034, 10005, "eng", $true, "superadmin", "now()", "now()"),
    @(10005, 110035, 10005, "eng", $true, "superadmin", "now()", "now()")
)

$startRow = 34
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

# Move / resize the active selection to reflect the post-edit UI state
# (mirrors selecting the row(s) right below the newly entered data, i.e.
# row 37 through the end of the sheet - the same as clicking the row 37
# header in Excel).
$ws.Rows("37:1048576").Select()
